$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 319, shifting existing rows (319-372) down to (320-373)
$ws.Rows.Item(319).Insert()

# Populate the new row 319 with the values from the diff
$ws.Cells.Item(319, 1).Value = 9
$ws.Cells.Item(319, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(319, 3).Value = "Metropolitana"
$ws.Cells.Item(319, 4).Value = 44995
$ws.Cells.Item(319, 5).Value = 13
$ws.Cells.Item(319, 6).Value = 100112030
$ws.Cells.Item(319, 7).Value = "Poroto granado"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 70
$ws.Cells.Item(319, 11).Value = 32000
$ws.Cells.Item(319, 12).Value = 34000
$ws.Cells.Item(319, 13).Value = 33000
$ws.Cells.Item(319, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(319, 15).Value = "Región Metropolitana"
$ws.Cells.Item(319, 16).Value = 1320
$ws.Cells.Item(319, 17).Value = 25
$ws.Cells.Item(319, 18).Value = "Hortaliza"
